# XF-965 AUTO_TC 6.2.2 Read overview layout of a Tenant
# Populate the "6_Tenants" parameter sheet with the Tenant overview labels
# and values used by the "Compare text from web elements" helper method.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("6_Tenants")

# ---- Row 1: column headers (bold) -----------------------------------
$ws.Range("B1").Value = 'Title Name'
$ws.Range("C1").Value = 'Title Location'
$ws.Range("D1").Value = 'Title URL'
$ws.Range("E1").Value = 'Title Description'
$ws.Range("F1").Value = 'Title User'
$ws.Range("G1").Value = 'Title Quick Search'
$ws.Range("H1").Value = 'Tenant Name'
$ws.Range("I1").Value = 'Tenant Location'
$ws.Range("J1").Value = 'Tenant URL'
$ws.Range("K1").Value = 'Tenant Description'
$ws.Range("L1").Value = 'Find User'
$ws.Range("M1").Value = 'Tenant Management Title 1'
$ws.Range("N1").Value = 'Tenant Management Title 2'
$ws.Range("O1").Value = 'Tenant Management Title 3'
$ws.Range("P1").Value = 'Tenant Management Table Title 1'
$ws.Range("Q1").Value = 'Tenant Management Table Title 2'
$ws.Range("R1").Value = 'Tenant Management Table Title 3'
$ws.Range("S1").Value = 'Tenant Management Table Title 4'
$ws.Range("T1").Value = 'Tenant Management Table Title 5'
$ws.Range("U1").Value = 'Tenant Creator Title'
$ws.Range("V1").Value = 'Web Site URL Title'
$ws.Range("W1").Value = 'No Of Tenants Title'
$ws.Range("X1").Value = 'No of roles title'
$ws.Range("Y1").Value = 'Description Title'
$ws.Range("Z1").Value = 'Tenant Location Title'
$ws.Range("AA1").Value = 'Date Created Title'
$ws.Range("AB1").Value = 'Overview Title'
$ws.Range("AC1").Value = 'Tenant Name Title'
$ws.Range("AD1").Value = 'Overview Table DescTitle'
$ws.Range("AE1").Value = 'Overview Table User Title'
$ws.Range("AF1").Value = 'Overview Table Date Title'
$ws.Range("B1:AF1").Font.Bold = $true

# ---- Row 2: sample tenant data ---------------------------------------
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = '1'

$ws.Range("B2:G2").NumberFormat = "@"
$ws.Range("B2").Value = 'Tenant Name'
$ws.Range("C2").Value = 'Tenant Location'
$ws.Range("D2").Value = 'Tenant URL'
$ws.Range("E2").Value = 'Tenant Description'
$ws.Range("F2").Value = 'Find User'
$ws.Range("G2").Value = 'Quick Search'

$ws.Range("H2").Value = 'QA tenant test'

$ws.Range("I2:AF2").NumberFormat = "@"
$ws.Range("I2").Value = 'QA Tenant Test '
$ws.Range("J2").Value = 'xform-stage.janeirodigital.com'
$ws.Range("K2").Value = 'This is an automation webdriver test'
$ws.Range("L2").Value = 'Gerardo Soto'
$ws.Range("M2").Value = 'Active Tenants'
$ws.Range("N2").Value = 'New Tenants in 30 days'
$ws.Range("O2").Value = 'New Users added'
$ws.Range("P2").Value = 'Tenant'
$ws.Range("Q2").Value = 'Parent'
$ws.Range("R2").Value = 'No. of users'
$ws.Range("S2").Value = 'Tenant Administrators'
$ws.Range("T2").Value = 'Actions'
$ws.Range("U2").Value = 'Tenant Creator'
$ws.Range("V2").Value = 'Website Url'
$ws.Range("W2").Value = 'No. of tenants'
$ws.Range("X2").Value = 'No. of roles'
$ws.Range("Y2").Value = 'Description'
$ws.Range("Z2").Value = 'Tenant Location'
$ws.Range("AA2").Value = 'Date Created'
$ws.Range("AB2").Value = 'Overview: '
$ws.Range("AC2").Value = 'Tenant Name'
$ws.Range("AD2").Value = 'Description'
$ws.Range("AE2").Value = 'User'
$ws.Range("AF2").Value = 'Date/Time'

# Trailing formatted-but-empty cells that round out the used range
$ws.Range("AG2:BL2").NumberFormat = "@"
$ws.Range("F3:R3").NumberFormat = "@"

# ---- Make "6_Tenants" the active sheet/tab, matching the saved view --
$ws.Activate()
[void]$ws.Range("L19").Select()
